$d = $word.ActiveDocument

# Remove every paragraph after the first one ({OverarchingTheme}).
# This drops the CompanyWorked/Country, JobScope/TimePeriodWorked and all
# five ExperiencePoints[...] paragraphs, leaving only the opening paragraph.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 2; $i--) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Delete()
}

# The removed paragraphs used to host the document's "_GoBack" bookmark
# (originally anchored around the first ExperiencePoints run). Recreate it
# at the end of the surviving paragraph, right after the {OverarchingTheme}
# run and before the paragraph mark.
$p1 = $d.Paragraphs.Item(1)
$endPos = $p1.Range.End - 1

# Work around degenerate-range placement quirks by temporarily appending a
# placeholder character, anchoring the zero-length bookmark immediately
# before it, then deleting the placeholder again.
$tmpRange = $d.Range($endPos, $endPos)
$tmpRange.InsertAfter("Z")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$zRange = $d.Range($endPos, $endPos + 1)
$zRange.Delete()
